# NYPD 111th Precinct CompStat weekly report refresh:
# bump the report Volume/Number and reporting week dates, and replace the
# crime-complaint figures (Week-to-Date / 28-Day / Year-to-Date counts and
# their derived percentage changes) with newly collected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 30   Number  13" -> "...  14"
# ---------------------------------------------------------------------
$hdr = $ws.Range("A8")
$hdrText = $hdr.Characters().Text
$hdr.Characters($hdrText.Length - 1, 2).Text = "14"

# ---------------------------------------------------------------------
# Header text: reporting week dates
# "Report Covering the Week  3/27/2023  Through  4/2/2023"
#   -> "...  4/3/2023  Through  4/9/2023"
# (replace the later date first so the earlier date's offset is untouched)
# ---------------------------------------------------------------------
$wk = $ws.Range("C9")
$wkText = $wk.Characters().Text
$throughPos = $wkText.IndexOf("4/2/2023") + 1
$wk.Characters($throughPos, 8).Text = "4/9/2023"
$startPos = $wkText.IndexOf("3/27/2023") + 1
$wk.Characters($startPos, 9).Text = "4/3/2023"

# ---------------------------------------------------------------------
# Crime complaints grid (rows 15-27). A handful of cells that previously
# held the "no data" placeholders ("0" / "***.*") now carry real numbers,
# so give them the same numeric styling used elsewhere in their column.
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 1
$ws.Range("D17").NumberFormat = $ws.Range("C17").NumberFormat
$ws.Range("E17").Value = 300
$ws.Range("E17").NumberFormat = $ws.Range("E16").NumberFormat

$ws.Range("D25").Value = 6
$ws.Range("D25").NumberFormat = $ws.Range("C25").NumberFormat
$ws.Range("E25").Value = -66.666666666666
$ws.Range("E25").NumberFormat = $ws.Range("E24").NumberFormat

$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = $ws.Range("G27").NumberFormat
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = $ws.Range("G27").NumberFormat

# Row 15
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 2
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0

# Row 16
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 20
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = -16.666666666666
$ws.Range("L16").Value = 566.666666666667
$ws.Range("M16").Value = -13.043478260869
$ws.Range("N16").Value = -81.132075471698

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 225
$ws.Range("I17").Value = 26
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 30
$ws.Range("L17").Value = 136.363636363636
$ws.Range("M17").Value = 62.5
$ws.Range("N17").Value = -13.333333333333

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -72.727272727272
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -29.629629629629
$ws.Range("I18").Value = 97
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = 5.434782608695
$ws.Range("L18").Value = 51.5625
$ws.Range("M18").Value = 32.876712328767
$ws.Range("N18").Value = -63.396226415094

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 8.333333333333
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = 6.122448979591
$ws.Range("I19").Value = 206
$ws.Range("J19").Value = 148
$ws.Range("K19").Value = 39.189189189189
$ws.Range("L19").Value = 139.53488372093
$ws.Range("M19").Value = 96.190476190476
$ws.Range("N19").Value = 50.364963503649

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 33.333333333333
$ws.Range("L20").Value = 122.222222222222
$ws.Range("M20").Value = 5.263157894736
$ws.Range("N20").Value = -95.321637426900

# Row 21 (TOTAL)
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -17.857142857142
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = 6.521739130434
$ws.Range("I21").Value = 392
$ws.Range("J21").Value = 314
$ws.Range("K21").Value = 24.840764331210
$ws.Range("L21").Value = 115.384615384615
$ws.Range("M21").Value = 52.529182879377
$ws.Range("N21").Value = -71.919770773639

# Row 24
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -60
$ws.Range("G24").Value = 62
$ws.Range("H24").Value = -45.161290322580
$ws.Range("I24").Value = 142
$ws.Range("J24").Value = 217
$ws.Range("K24").Value = -34.562211981566
$ws.Range("L24").Value = -0.699300699300
$ws.Range("M24").Value = 3.649635036496

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("F25").Value = 21
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 57
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = -5
$ws.Range("L25").Value = 83.870967741935
$ws.Range("M25").Value = 67.647058823529

# Row 26
$ws.Range("F26").Value = 3
$ws.Range("I26").Value = 3
$ws.Range("L26").Value = 200

# Row 27
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 5
$ws.Range("K27").Value = -28.571428571428
$ws.Range("L27").Value = 66.666666666666
